$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new tracked problem: "Single Number" (LeetCode 136) as row 11 ---
# (Inserted first so its new shared strings land before the Phase notes below,
#  matching the shared-strings append order.)
$ws.Range("A11").Value = 136
$ws.Range("B11").Value = "Single Number"
$ws.Range("C11").Value = "Easy"
$ws.Range("D11").Value = "Bit Manupilation"
$ws.Range("E11").Value = "Accepted"
$ws.Range("F11").Value = "O(n)"
$ws.Range("G11").Value = "O(1)"
$ws.Range("H11").Value = "Phase 2-3"

# --- Extend the phase-notes list (columns K/L) with Phase 4 - Phase 6 ---
$ws.Range("K6").Value = "Phase 4"
$ws.Range("L6").Value = "Go back to Step 2"

$ws.Range("K7").Value = "Phase 5"
$ws.Range("L7").Value = "Completion of all phases"

$ws.Range("K8").Value = "Phase 6"
$ws.Range("L8").Value = "Revision"

# --- Update the "Delete Node in Linked List" row's phase note ---
$ws.Range("H10").Value = "Phase 3-4"

# --- Move the active selection to the newly added solution's notes cell ---
$ws.Range("H11").Select() | Out-Null
